$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "en"
$ws2 = $wb.Worksheets.Item(2)   # "de"

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) New rows 212-215 on both sheets - set VALUES first, in the precise order
#    needed so that newly created shared-strings are appended in the same
#    sequence as the target workbook (Level, Niveau, Referenz, Courses,
#    Kurse, Course, Kurs).
# ---------------------------------------------------------------------------

# Row 212: Level / Niveau
$ws1.Cells.Item(212, 1).Value2 = "Level"
$ws1.Cells.Item(212, 2).Value2 = "Level"
$ws2.Cells.Item(212, 1).Value2 = "Level"
$ws2.Cells.Item(212, 2).Value2 = "Niveau"

# Row 206 (existing row): B206 value change - "Referenz" is a brand new string
$ws2.Cells.Item(206, 2).Value2 = "Referenz"
$ws1.Cells.Item(206, 2).Value2 = "Reference"

# Row 215: Courses / Kurse
$ws1.Cells.Item(215, 1).Value2 = "Courses"
$ws1.Cells.Item(215, 2).Value2 = "Courses"
$ws2.Cells.Item(215, 1).Value2 = "Courses"
$ws2.Cells.Item(215, 2).Value2 = "Kurse"

# Row 214: Course / Kurs
$ws1.Cells.Item(214, 1).Value2 = "Course"
$ws1.Cells.Item(214, 2).Value2 = "Course"
$ws2.Cells.Item(214, 1).Value2 = "Course"
$ws2.Cells.Item(214, 2).Value2 = "Kurs"

# Row 213: References / Referenzen (both strings already exist elsewhere)
$ws1.Cells.Item(213, 1).Value2 = "References"
$ws1.Cells.Item(213, 2).Value2 = "References"
$ws2.Cells.Item(213, 1).Value2 = "References"
$ws2.Cells.Item(213, 2).Value2 = "Referenzen"

# ---------------------------------------------------------------------------
# 2) Apply cell formatting to the new rows by copying the number/alignment
#    format from nearby cells that already carry the desired style, using
#    Copy + PasteSpecial(Formats) so the existing style table is reused
#    instead of synthesizing brand-new style entries.
#       A208 -> style "1" (wrap + vertical-center)
#       A209 -> style "0" (default / no explicit style)
#       B209 -> style "5" (number format + wrap)
# ---------------------------------------------------------------------------

# Row 212: A=no style, B=style 5
$ws1.Range("A209").Copy()
$ws1.Range("A212").PasteSpecial($xlPasteFormats)
$ws1.Range("B209").Copy()
$ws1.Range("B212").PasteSpecial($xlPasteFormats)
$ws2.Range("A209").Copy()
$ws2.Range("A212").PasteSpecial($xlPasteFormats)
$ws2.Range("B209").Copy()
$ws2.Range("B212").PasteSpecial($xlPasteFormats)

# Row 213: A=style 1, B=style 1 (en) / style 5 (de)
$ws1.Range("A208").Copy()
$ws1.Range("A213").PasteSpecial($xlPasteFormats)
$ws1.Range("A208").Copy()
$ws1.Range("B213").PasteSpecial($xlPasteFormats)
$ws2.Range("A208").Copy()
$ws2.Range("A213").PasteSpecial($xlPasteFormats)
$ws2.Range("B209").Copy()
$ws2.Range("B213").PasteSpecial($xlPasteFormats)

# Row 214: A=no style, B=no style (en) / style 5 (de)
$ws1.Range("A209").Copy()
$ws1.Range("A214").PasteSpecial($xlPasteFormats)
$ws1.Range("A209").Copy()
$ws1.Range("B214").PasteSpecial($xlPasteFormats)
$ws2.Range("A209").Copy()
$ws2.Range("A214").PasteSpecial($xlPasteFormats)
$ws2.Range("B209").Copy()
$ws2.Range("B214").PasteSpecial($xlPasteFormats)

# Row 215: A=no style, B=no style (en) / style 5 (de)
$ws1.Range("A209").Copy()
$ws1.Range("A215").PasteSpecial($xlPasteFormats)
$ws1.Range("A209").Copy()
$ws1.Range("B215").PasteSpecial($xlPasteFormats)
$ws2.Range("A209").Copy()
$ws2.Range("A215").PasteSpecial($xlPasteFormats)
$ws2.Range("B209").Copy()
$ws2.Range("B215").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Update selections to match the recorded view state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B214:B215").Select()

$ws2.Activate()
$ws2.Range("A214:B215").Select()

$ws1.Activate()

Write-Output "done"
